$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target header row (A1:E1) - columns shifted one left, MODEL_CONDITION fixed to MODELCONDITION
$ws.Range("A1").Value = "EL_Astral15"
$ws.Range("B1").Value = "FNRATE_EXACT_ASTRAL"
$ws.Range("C1").Value = "TAXON"
$ws.Range("D1").Value = "MODELCONDITION"
$ws.Range("E1").Value = "GENE"

# Target data rows (A2:E5) - columns shifted one left (old column A dropped)
$ws.Range("A2").Value = 51
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "11-texon"
$ws.Range("D2").Value = "estimated_15genes_weakILS"
$ws.Range("E2").Value = 0

$ws.Range("A3").Value = 58
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "11-texon"
$ws.Range("D3").Value = "estimated_15genes_weakILS"
$ws.Range("E3").Value = 1

$ws.Range("A4").Value = 51
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "11-texon"
$ws.Range("D4").Value = "estimated_15genes_weakILS"
$ws.Range("E4").Value = 9

$ws.Range("A5").Value = 58
$ws.Range("B5").Value = 0.125
$ws.Range("C5").Value = "11-texon"
$ws.Range("D5").Value = "estimated_15genes_weakILS"
$ws.Range("E5").Value = 12

# Remove the now-unused old column F entirely
$ws.Range("F1:F5").Clear()

# Data cells in column A no longer carry the header-style formatting
$ws.Range("A2:A5").ClearFormats()

# A1 (previously blank/unstyled) now holds a header label and must match the
# bold/centered/bordered style used by the rest of the header row
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

